$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update calibration data values (row 4)
$ws.Range("C4").Value = 1430
$ws.Range("F4").Value = 1950

# Update calibration data values (row 5)
$ws.Range("C5").Value = 1580
$ws.Range("D5").Value = 1570
$ws.Range("G5").Value = 2000

# Update calibration data values (row 6)
$ws.Range("C6").Value = 1400
$ws.Range("D6").Value = 1680
$ws.Range("G6").Value = 2150

# Update calibration data values (row 8)
$ws.Range("D8").Value = 1200
$ws.Range("G8").Value = 730

# Update calibration data values (row 9)
$ws.Range("C9").Value = 1430
$ws.Range("D9").Value = 1230
$ws.Range("E9").Value = 1830
$ws.Range("F9").Value = 900

# Update the active selection on the sheet
$ws.Range("G11").Select()
